$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 328.16
$ws.Range("I15").Value = 328.16
$ws.Range("K15").Value = 984.48
$ws.Range("M15").Value = -815.48

$ws.Range("H17").Value = 993.17645
$ws.Range("I17").Value = 652.6667
$ws.Range("J17").Value = 1262
$ws.Range("K17").Value = 1958.0001
$ws.Range("L17").Value = 3786
$ws.Range("M17").Value = -1790.0001
$ws.Range("N17").Value = -4122

$ws.Range("H33").Value = 246.17857
$ws.Range("I33").Value = 245.92308
$ws.Range("K33").Value = 245.92308
$ws.Range("M33").Value = -16.92308

$ws.Range("H62").Value = 22009.889
$ws.Range("I62").Value = 13278
$ws.Range("J62").Value = 32924.75
$ws.Range("K62").Value = 13278
$ws.Range("L62").Value = 32924.75
$ws.Range("M62").Value = -12654
$ws.Range("N62").Value = -34172.75

$ws.Range("H65").Value = 22009.889
$ws.Range("I65").Value = 13278
$ws.Range("J65").Value = 32924.75
$ws.Range("K65").Value = 66390
$ws.Range("L65").Value = 164623.75
$ws.Range("M65").Value = -63270
$ws.Range("N65").Value = -170863.75

$ws.Range("H112").Value = 5338.237
$ws.Range("J112").Value = 5922
$ws.Range("L112").Value = 17766
$ws.Range("N112").Value = -19982

$ws.Range("H138").Value = 2826.65
$ws.Range("I138").Value = 1489
$ws.Range("J138").Value = 3081.4404
$ws.Range("K138").Value = 4467
$ws.Range("L138").Value = 9244.3212
$ws.Range("M138").Value = 673
$ws.Range("N138").Value = -19524.3212

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 90028
$ws.Range("J34").Value = 90028
$ws.Range("L34").Value = 90028
$ws.Range("N34").Value = -90570

$ws.Range("H45").Value = 2132.95
$ws.Range("I45").Value = 1842.1818
$ws.Range("J45").Value = 2488.3333
$ws.Range("K45").Value = 1842.1818
$ws.Range("L45").Value = 2488.3333
$ws.Range("M45").Value = -1465.1818
$ws.Range("N45").Value = -3242.3333

$ws.Range("H110").Value = 82618.17999999999
$ws.Range("I110").Value = 129271.43
$ws.Range("J110").Value = 975
$ws.Range("K110").Value = 129271.43
$ws.Range("L110").Value = 975
$ws.Range("M110").Value = -127226.43
$ws.Range("N110").Value = -5065

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5340.0884
$ws.Range("I31").Value = 1790.3914
$ws.Range("J31").Value = 7154.378
$ws.Range("K31").Value = 1790.3914
$ws.Range("L31").Value = 7154.378
$ws.Range("M31").Value = -1495.3914
$ws.Range("N31").Value = -7744.378

$ws.Range("H34").Value = 5340.0884
$ws.Range("I34").Value = 1790.3914
$ws.Range("J34").Value = 7154.378
$ws.Range("K34").Value = 1790.3914
$ws.Range("L34").Value = 7154.378
$ws.Range("M34").Value = -1588.3914
$ws.Range("N34").Value = -7558.378

$ws.Range("H62").Value = 4330.727
$ws.Range("I62").Value = 4405.853
$ws.Range("J62").Value = 3763.111
$ws.Range("K62").Value = 4405.853
$ws.Range("L62").Value = 3763.111
$ws.Range("M62").Value = -3781.853
$ws.Range("N62").Value = -5011.111

$ws.Range("H65").Value = 4330.727
$ws.Range("I65").Value = 4405.853
$ws.Range("J65").Value = 3763.111
$ws.Range("K65").Value = 22029.265
$ws.Range("L65").Value = 18815.555
$ws.Range("M65").Value = -18909.265
$ws.Range("N65").Value = -25055.555

$ws.Range("H132").Value = 19233190
$ws.Range("I132").Value = 29413830
$ws.Range("J132").Value = 3091.2222
$ws.Range("K132").Value = 88241490
$ws.Range("L132").Value = 9273.6666
$ws.Range("M132").Value = -88238960
$ws.Range("N132").Value = -14333.6666

$ws.Range("H134").Value = 6544.2104
$ws.Range("I134").Value = 7782
$ws.Range("J134").Value = 1902.5
$ws.Range("K134").Value = 23346
$ws.Range("L134").Value = 5707.5
$ws.Range("M134").Value = -20811
$ws.Range("N134").Value = -10777.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 676.6613
$ws.Range("I5").Value = 473.07144
$ws.Range("J5").Value = 1104.2
$ws.Range("K5").Value = 1419.21432
$ws.Range("L5").Value = 3312.6
$ws.Range("M5").Value = -1307.21432
$ws.Range("N5").Value = -3536.6

$ws.Range("H122").Value = 3267
$ws.Range("I122").Value = 484.22223
$ws.Range("J122").Value = 5444.826
$ws.Range("K122").Value = 4358.00007
$ws.Range("L122").Value = 49003.434
$ws.Range("M122").Value = -1908.00007
$ws.Range("N122").Value = -53903.434

$ws.Range("H135").Value = 676.6613
$ws.Range("I135").Value = 473.07144
$ws.Range("J135").Value = 1104.2
$ws.Range("K135").Value = 4257.64296
$ws.Range("L135").Value = 9937.800000000001
$ws.Range("M135").Value = -1722.64296
$ws.Range("N135").Value = -15007.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 75201.53
$ws.Range("I113").Value = 93443.664
$ws.Range("K113").Value = 93443.664
$ws.Range("M113").Value = -91273.664

$ws.Range("H136").Value = 21096.475
$ws.Range("I136").Value = 44500
$ws.Range("J136").Value = 18343.117
$ws.Range("K136").Value = 133500
$ws.Range("L136").Value = 55029.351
$ws.Range("M136").Value = -130950
$ws.Range("N136").Value = -60129.351

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5968.4
$ws.Range("I7").Value = 5742.6665
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 5742.6665
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -5630.6665
$ws.Range("N7").Value = -8224

$ws.Range("H68").Value = 2550.3333
$ws.Range("I68").Value = 2001
$ws.Range("J68").Value = 2825
$ws.Range("K68").Value = 2001
$ws.Range("L68").Value = 2825
$ws.Range("M68").Value = -1252
$ws.Range("N68").Value = -4323

$ws.Range("H71").Value = 2550.3333
$ws.Range("I71").Value = 2001
$ws.Range("J71").Value = 2825
$ws.Range("K71").Value = 10005
$ws.Range("L71").Value = 14125
$ws.Range("M71").Value = -6261
$ws.Range("N71").Value = -21613

$ws.Range("H100").Value = 49842.844
$ws.Range("I100").Value = 55193.53
$ws.Range("K100").Value = 55193.53
$ws.Range("M100").Value = -54652.53

$ws.Range("H126").Value = 5968.4
$ws.Range("I126").Value = 5742.6665
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 17227.9995
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -14757.9995
$ws.Range("N126").Value = -28940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2137
$ws.Range("I122").Value = 2200.3845
$ws.Range("J122").Value = 1725
$ws.Range("K122").Value = 6601.1535
$ws.Range("L122").Value = 5175
$ws.Range("M122").Value = -4151.1535
$ws.Range("N122").Value = -10075
